$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.386.39"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.780.36"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.81"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5354"
$ws.Range("E7").Value = "  +12.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3784"
$ws.Range("E8").Value = "  +9.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.86"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07408"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.097"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.105"
$ws.Range("E14").Value = "  +4.79%  "
$ws.Range("D15").Value = "1.776.22"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.989"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.61"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001056"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06440"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.81"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.905"
$ws.Range("D23").Value = "27.424.67"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  +4.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.080"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.82"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.367"
$ws.Range("E28").Value = "  +13.78%  "
$ws.Range("D29").Value = "1.981.35"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.07"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.081"
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1033"
$ws.Range("E32").Value = "  +13.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.594"
$ws.Range("E33").Value = "  +5.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.626"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02256"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2059"
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.920"
$ws.Range("E38").Value = "  +4.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.28"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.270"
$ws.Range("E40").Value = "  +11.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6117"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.32"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5773"
$ws.Range("E45").Value = "  +3.64%  "
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.01"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.124"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.86"
$ws.Range("E51").Value = "  +2.76%  "
